$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.077679419419598
$ws.Range("D2").Value = 1.077051058944553
$ws.Range("E2").Value = 1.0804264580099
$ws.Range("F2").Value = 1.089170678772191
$ws.Range("I2").Value = 1.050287125591677
$ws.Range("J2").Value = 1.082572676276497
$ws.Range("K2").Value = 1.079733377003637
$ws.Range("L2").Value = 1.083099922956953
$ws.Range("M2").Value = 1.0918214844453
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.079274077384807
$ws.Range("D3").Value = 1.078302463502724
$ws.Range("E3").Value = 1.081814909288512
$ws.Range("F3").Value = 1.090561679395764
$ws.Range("I3").Value = 1.050673018726119
$ws.Range("J3").Value = 1.083824640796738
$ws.Range("K3").Value = 1.080801342285464
$ws.Range("L3").Value = 1.084305246673492
$ws.Range("M3").Value = 1.093031003619709
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.080304654526588
$ws.Range("D4").Value = 1.079110843163894
$ws.Range("E4").Value = 1.08271230434714
$ws.Range("F4").Value = 1.091460663694052
$ws.Range("I4").Value = 1.050920746271911
$ws.Range("J4").Value = 1.084633051092644
$ws.Range("K4").Value = 1.081490460681232
$ws.Range("L4").Value = 1.085083608082098
$ws.Range("M4").Value = 1.093812017004684
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.080737612399024
$ws.Range("D5").Value = 1.07945036516801
$ws.Range("E5").Value = 1.083089329772964
$ws.Range("F5").Value = 1.091838342704078
$ws.Range("I5").Value = 1.05102442116833
$ws.Range("J5").Value = 1.08497250717956
$ws.Range("K5").Value = 1.081779709941692
$ws.Range("L5").Value = 1.08541046253944
$ws.Range("M5").Value = 1.094139971259656
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.080810290761981
$ws.Range("D6").Value = 1.079507353751882
$ws.Range("E6").Value = 1.083152620136098
$ws.Range("F6").Value = 1.091901741952299
$ws.Range("I6").Value = 1.05104180115396
$ws.Range("J6").Value = 1.085029480092796
$ws.Range("K6").Value = 1.08182824956646
$ws.Range("L6").Value = 1.085465321385388
$ws.Range("M6").Value = 1.094195013879437
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.080310440886963
$ws.Range("D7").Value = 1.079115381127348
$ws.Range("E7").Value = 1.082717343112671
$ws.Range("F7").Value = 1.091465711250011
$ws.Range("I7").Value = 1.050922133424181
$ws.Range("J7").Value = 1.084637588484167
$ws.Range("K7").Value = 1.081494327429913
$ws.Range("L7").Value = 1.085087976969018
$ws.Range("M7").Value = 1.093816400646543
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.078218610023318
$ws.Range("D8").Value = 1.077474262429914
$ws.Range("E8").Value = 1.080895907699891
$ws.Range("F8").Value = 1.089641002472957
$ws.Range("I8").Value = 1.050417949399115
$ws.Range("J8").Value = 1.082996138053326
$ws.Range("K8").Value = 1.080094702859812
$ws.Range("L8").Value = 1.083507594925812
$ws.Range("M8").Value = 1.092230587244521
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.074522424224321
$ws.Range("D9").Value = 1.074571719545457
$ws.Range("E9").Value = 1.077678174874104
$ws.Range("F9").Value = 1.086417042954855
$ws.Range("I9").Value = 1.049514320603024
$ws.Range("J9").Value = 1.080090442810736
$ws.Range("K9").Value = 1.077613400672973
$ws.Range("L9").Value = 1.080710529878791
$ws.Range("M9").Value = 1.089423473702824
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.072050994147459
$ws.Range("D10").Value = 1.072629164192416
$ws.Range("E10").Value = 1.075527166485464
$ws.Range("F10").Value = 1.084261580086386
$ws.Range("I10").Value = 1.048901551965604
$ws.Range("J10").Value = 1.078144021445738
$ws.Range("K10").Value = 1.075948800878152
$ws.Range("L10").Value = 1.078837239361449
$ws.Range("M10").Value = 1.087543166998733
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.070978991348959
$ws.Range("D11").Value = 1.071786153617815
$ws.Range("E11").Value = 1.074594281183619
$ws.Range("F11").Value = 1.083326695066479
$ws.Range("I11").Value = 1.04863373117173
$ws.Range("J11").Value = 1.077298909981148
$ws.Range("K11").Value = 1.075225471857247
$ws.Range("L11").Value = 1.078023968494263
$ws.Range("M11").Value = 1.086726783668431
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.070580512906882
$ws.Range("D12").Value = 1.071472734056824
$ws.Range("E12").Value = 1.07424753547534
$ws.Range("F12").Value = 1.082979196010138
$ws.Range("I12").Value = 1.048533874133648
$ws.Range("J12").Value = 1.076984645206497
$ws.Range("K12").Value = 1.074956406354372
$ws.Range("L12").Value = 1.0777215574259
$ws.Range("M12").Value = 1.086423205383225
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.07066600111616
$ws.Range("D13").Value = 1.071539976753689
$ws.Range("E13").Value = 1.074321924162378
$ws.Range("F13").Value = 1.08305374676951
$ws.Range("I13").Value = 1.048555310886747
$ws.Range("J13").Value = 1.07705207215066
$ws.Range("K13").Value = 1.075014139518919
$ws.Range("L13").Value = 1.077786440494267
$ws.Range("M13").Value = 1.086488339325179
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.070946058946694
$ws.Range("D14").Value = 1.071760252171461
$ws.Range("E14").Value = 1.074565623822622
$ws.Range("F14").Value = 1.083297975648484
$ws.Range("I14").Value = 1.048625484653316
$ws.Range("J14").Value = 1.077272939996231
$ws.Range("K14").Value = 1.075203238789111
$ws.Range("L14").Value = 1.077998977793133
$ws.Range("M14").Value = 1.086701696711021
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.0711185732453
$ws.Range("D15").Value = 1.071895932773901
$ws.Range("E15").Value = 1.074715744446884
$ws.Range("F15").Value = 1.083448420944708
$ws.Range("I15").Value = 1.048668671060174
$ws.Range("J15").Value = 1.077408976998335
$ws.Range("K15").Value = 1.07531969736718
$ws.Range("L15").Value = 1.078129885621675
$ws.Range("M15").Value = 1.086833108352179
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.07212209947403
$ws.Range("D16").Value = 1.072685072026312
$ws.Range("E16").Value = 1.075589047091703
$ws.Range("F16").Value = 1.084323591921081
$ws.Range("I16").Value = 1.048919273702947
$ws.Range("J16").Value = 1.078200059653883
$ws.Range("K16").Value = 1.075996751654728
$ws.Range("L16").Value = 1.078891168130565
$ws.Range("M16").Value = 1.08759730077002
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.07275108079777
$ws.Range("D17").Value = 1.073179572478614
$ws.Range("E17").Value = 1.076136444023509
$ws.Range("F17").Value = 1.08487214188839
$ws.Range("I17").Value = 1.049075802139898
$ws.Range("J17").Value = 1.07869566435924
$ws.Range("K17").Value = 1.076420763736316
$ws.Range("L17").Value = 1.079368127036672
$ws.Range("M17").Value = 1.088076065113103
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.073117776431832
$ws.Range("D18").Value = 1.073467826133974
$ws.Range("E18").Value = 1.07645558867553
$ws.Range("F18").Value = 1.085191952239887
$ws.Range("I18").Value = 1.049166862662682
$ws.Range("J18").Value = 1.078984520967244
$ws.Range("K18").Value = 1.076667837557831
$ws.Range("L18").Value = 1.079646124665479
$ws.Range("M18").Value = 1.088355108638812
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.073242780181227
$ws.Range("D19").Value = 1.073566082837924
$ws.Range("E19").Value = 1.076564384774239
$ws.Range("F19").Value = 1.085300974160931
$ws.Range("I19").Value = 1.049197871343088
$ws.Range("J19").Value = 1.079082976306127
$ws.Range("K19").Value = 1.076752041937877
$ws.Range("L19").Value = 1.079740880208169
$ws.Range("M19").Value = 1.08845021958468
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.072683615591902
$ws.Range("D20").Value = 1.073126535936783
$ws.Range("E20").Value = 1.076077728329846
$ws.Range("F20").Value = 1.084813303192475
$ws.Range("I20").Value = 1.049059032949968
$ws.Range("J20").Value = 1.078642513601924
$ws.Range("K20").Value = 1.076375296666316
$ws.Range("L20").Value = 1.079316975055972
$ws.Range("M20").Value = 1.08802472016867
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.070863596949504
$ws.Range("D21").Value = 1.071695394553104
$ws.Range("E21").Value = 1.074493866772488
$ws.Range("F21").Value = 1.083226063047529
$ws.Range("I21").Value = 1.048604830642213
$ws.Range("J21").Value = 1.077207909718579
$ws.Range("K21").Value = 1.075147564553889
$ws.Range("L21").Value = 1.077936399913184
$ws.Range("M21").Value = 1.086638877642016
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.069717602744563
$ws.Range("D22").Value = 1.070793910084483
$ws.Range("E22").Value = 1.073496692359458
$ws.Range("F22").Value = 1.082226703307469
$ws.Range("I22").Value = 1.048317076135127
$ws.Range("J22").Value = 1.076303872682771
$ws.Range("K22").Value = 1.074373387423526
$ws.Range("L22").Value = 1.077066487555471
$ws.Range("M22").Value = 1.085765589229425
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.070325278146968
$ws.Range("D23").Value = 1.071271964537073
$ws.Range("E23").Value = 1.074025442723107
$ws.Range("F23").Value = 1.08275661791412
$ws.Range("I23").Value = 1.048469827736632
$ws.Range("J23").Value = 1.076783316137578
$ws.Range("K23").Value = 1.07478400907725
$ws.Range("L23").Value = 1.077527826051142
$ws.Range("M23").Value = 1.086228723501843
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.072714100766296
$ws.Range("D24").Value = 1.07315050142384
$ws.Range("E24").Value = 1.076104259865677
$ws.Range("F24").Value = 1.084839890328896
$ws.Range("I24").Value = 1.049066610966351
$ws.Range("J24").Value = 1.078666530825121
$ws.Range("K24").Value = 1.076395842036168
$ws.Range("L24").Value = 1.079340089066146
$ws.Range("M24").Value = 1.088047921393058
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.075479227241966
$ws.Range("D25").Value = 1.075323397667219
$ws.Range("E25").Value = 1.078511039105214
$ws.Range("F25").Value = 1.087251570768124
$ws.Range("I25").Value = 1.049749744496203
$ws.Range("J25").Value = 1.080843244094679
$ws.Range("K25").Value = 1.078256685691701
$ws.Range("L25").Value = 1.081435124226183
$ws.Range("M25").Value = 1.090150722027586
